# Update the "as of" date in the confidential disclosure note (shared string).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet carries legacy sheet protection (no known password); unprotect
# so the cell values below can be written, then re-protect at the end.
$ws.Unprotect()

$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-27 for illustrative purposes only and are subject to change."
$ws.Rows(80).AutoFit()

# Refresh the Weight (D) and Percent Change (E) columns for each holding row.
$ws.Range("D2").Value = 0.0676855416704667
$ws.Range("E2").Value = -0.002449524940617676
$ws.Range("D3").Value = 0.04077946703473535
$ws.Range("E3").Value = 0.002472865943091707
$ws.Range("D4").Value = 0.03441612368210189
$ws.Range("E4").Value = 0.001605811508315957
$ws.Range("D5").Value = 0.03030934784528904
$ws.Range("E5").Value = -0.001866805591734022
$ws.Range("D6").Value = 0.02763206637945034
$ws.Range("E6").Value = -0.00820371180079027
$ws.Range("D7").Value = 0.02401391045956685
$ws.Range("E7").Value = 0.004582890541976692
$ws.Range("D8").Value = 0.1756722977819319
$ws.Range("E8").Value = 0.02218741842860861
$ws.Range("D9").Value = 0.02421344155670366
$ws.Range("E9").Value = -0.005727516451377146
$ws.Range("D10").Value = 0.02198239360566473
$ws.Range("E10").Value = -0.0006094773731524228
$ws.Range("D11").Value = 0.02267046487109095
$ws.Range("E11").Value = -0.001076426264800867
$ws.Range("D12").Value = 0.02057339463515587
$ws.Range("E12").Value = 0.00200792315623799
$ws.Range("D13").Value = 0.01917165279088181
$ws.Range("E13").Value = 0.01192590713017028
$ws.Range("D14").Value = 0.01709061203157523
$ws.Range("E14").Value = -0.005228505034856612
$ws.Range("D15").Value = 0.01782908443744528
$ws.Range("E15").Value = -0.004133007702423308
$ws.Range("D16").Value = 0.01617047222391555
$ws.Range("E16").Value = -0.005839185670322822
$ws.Range("D17").Value = 0.0142975754148596
$ws.Range("E17").Value = 0.003408019723007749
$ws.Range("D18").Value = 0.01429095627774435
$ws.Range("E18").Value = -0.003278459821428825
$ws.Range("D19").Value = 0.01329186531654133
$ws.Range("E19").Value = 0.001748944033791
$ws.Range("D20").Value = 0.01243313196163759
$ws.Range("E20").Value = 0.01311063218390807
$ws.Range("D21").Value = 0.01220190078012353
$ws.Range("E21").Value = -0.0045292785506309
$ws.Range("D22").Value = 0.01288064145467644
$ws.Range("E22").Value = 0.0004488733279468526
$ws.Range("D23").Value = 0.01236008220745001
$ws.Range("E23").Value = 0.004129351949828219
$ws.Range("D24").Value = 0.01262771864309795
$ws.Range("E24").Value = -0.003435558025564611
$ws.Range("D25").Value = 0.01148344528287911
$ws.Range("E25").Value = 0.001875058595581125
$ws.Range("D26").Value = 0.009807168745602063
$ws.Range("E26").Value = 0.01844262295081944
$ws.Range("D27").Value = 0.01025838655039825
$ws.Range("E27").Value = 0.002658706096366537
$ws.Range("D28").Value = 0.01065465754227388
$ws.Range("E28").Value = -0.03386151508573909
$ws.Range("D29").Value = 0.01072990038405988
$ws.Range("E29").Value = -0.002590182576283961
$ws.Range("D30").Value = 0.01021101585827225
$ws.Range("E30").Value = 0.001249609497032234
$ws.Range("D31").Value = 0.009137758661922452
$ws.Range("E31").Value = -0.01344452008168817
$ws.Range("D32").Value = 0.00996347607988989
$ws.Range("E32").Value = -0.0003681885125184081
$ws.Range("D33").Value = 0.009308819494599694
$ws.Range("E33").Value = -0.01088865472427114
$ws.Range("D34").Value = 0.008986555240590961
$ws.Range("E34").Value = -0.001490868430860881
$ws.Range("D35").Value = 0.009079382657484345
$ws.Range("E35").Value = 0.007290294246816043
$ws.Range("D36").Value = 0.008798029455766253
$ws.Range("E36").Value = -0.007178985152553552
$ws.Range("D37").Value = 0.008575171881446183
$ws.Range("E37").Value = 0.00850945107065626
$ws.Range("D38").Value = 0.008830566900862905
$ws.Range("E38").Value = -0.04532646979138444
$ws.Range("D39").Value = 0.008637096700361743
$ws.Range("E39").Value = -0.005946225439503405
$ws.Range("D40").Value = 0.007691277830640493
$ws.Range("E40").Value = 0.01162332545311262
$ws.Range("D41").Value = 0.007525879151399188
$ws.Range("E41").Value = -0.01165624668856624
$ws.Range("D42").Value = 0.007757309704513106
$ws.Range("E42").Value = -0.006836499712147504
$ws.Range("D43").Value = 0.007723655778457257
$ws.Range("E43").Value = 0.009664429530201302
$ws.Range("D44").Value = 0.007346444711527956
$ws.Range("E44").Value = 0.003875379939209944
$ws.Range("D45").Value = 0.007879803615465201
$ws.Range("E45").Value = -0.007732167435835158
$ws.Range("D46").Value = 0.007422365416753473
$ws.Range("E46").Value = 0.01211964930376475
$ws.Range("D47").Value = 0.007419494465715533
$ws.Range("E47").Value = -0.00219269960015489
$ws.Range("D48").Value = 0.007124105503367515
$ws.Range("E48").Value = 0.00720906282183309
$ws.Range("D49").Value = 0.006933346756624408
$ws.Range("E49").Value = 0.001984126984126977
$ws.Range("D50").Value = 0.006645055423231296
$ws.Range("E50").Value = -0.008448844884488382
$ws.Range("D51").Value = 0.006435396248821754
$ws.Range("E51").Value = -0.0009541984732823749
$ws.Range("D52").Value = 0.006525990703796741
$ws.Range("E52").Value = 0.00267621468374224
$ws.Range("D53").Value = 0.005575626161598752
$ws.Range("E53").Value = -0.005792748337266618
$ws.Range("D54").Value = 0.005972216148034146
$ws.Range("E54").Value = -0.02675996154257021
$ws.Range("D55").Value = 0.006005670702490138
$ws.Range("E55").Value = 0.001693058460312757
$ws.Range("D56").Value = 0.005695631914984632
$ws.Range("E56").Value = -0.0003916981005267495
$ws.Range("D57").Value = 0.005608243355335345
$ws.Range("E57").Value = 0.1042034013992377
$ws.Range("D58").Value = 0.005584956752472056
$ws.Range("E58").Value = -0.009938313913639418
$ws.Range("D59").Value = 0.005182066623481184
$ws.Range("E59").Value = 0.01094182825484746
$ws.Range("D60").Value = 0.004904222361920573
$ws.Range("E60").Value = -0.008846103811630091
$ws.Range("D61").Value = 0.004602852251576845
$ws.Range("E61").Value = -0.0002858776443681599
$ws.Range("D62").Value = 0.004611863847890378
$ws.Range("E62").Value = 0.001037523776586502
$ws.Range("D63").Value = 0.004163516994132127
$ws.Range("E63").Value = -0.008197977321483418
$ws.Range("D64").Value = 0.004203710308663283
$ws.Range("E64").Value = 0.001062376688419908
$ws.Range("D65").Value = 0.004106257470653218
$ws.Range("E65").Value = 0.002952029520295163
$ws.Range("D66").Value = 0.003729843890123345
$ws.Range("E66").Value = -0.0161642078255293
$ws.Range("D67").Value = 0.003871477474661704
$ws.Range("E67").Value = 0.005870720553701858
$ws.Range("D68").Value = 0.003298722742592731
$ws.Range("E68").Value = 0.04260951552074266
$ws.Range("D69").Value = 0.003557905822406733
$ws.Range("E69").Value = 0.005962254000986267
$ws.Range("D70").Value = 0.003058041347245452
$ws.Range("E70").Value = 0.02983362019506597
$ws.Range("D71").Value = 0.003250395066787413
$ws.Range("E71").Value = -0.01067275136169588
$ws.Range("D72").Value = 0.002468698898068383
$ws.Range("E72").Value = -0.006218503682646426
$ws.Range("D73").Value = 0.002034786548139772
$ws.Range("E73").Value = -0.009308250048990829
$ws.Range("D74").Value = 0.00205631868092432
$ws.Range("E74").Value = 0.003141361256544517
$ws.Range("D75").Value = 0.001537952521296316
$ws.Range("E75").Value = 0.009437386569872919
$ws.Range("D76").Value = 0.001433322305691402
$ws.Range("E76").Value = -0.005396984365437008
$ws.Range("E77").Value = 0.003522539294308125

$ws.Protect()
